{"js": "// Update the date paragraph and the 25 division-answer table cells\n// (in document order) to the new values from the target revision.\n\nconst body = context.document.body;\n\n// --- 1. Update the date paragraph (first paragraph in the body) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text.trim() === \"2025-12-18 Thursday\") {\n  dateParagraph.insertText(\"2025-12-19 Friday\", Word.InsertLocation.replace);\n}\n\n// --- 2. Update the table cells, in row/column (document) order ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// New values, listed in the same left-to-right, top-to-bottom order as\n// they appear in the document. Only every 4th row (0, 4, 8, 12, 16) of\n// this 20-row table actually holds text; the rows in between are blank\n// spacer rows that must stay untouched.\nconst newValues = [\n  [\"88\u00f78=11, 0\", \"76\u00f79=8, 4\", \"73\u00f77=10, 3\", \"50\u00f76=8, 2\", \"67\u00f77=9, 4\"],\n  [\"96\u00f77=13, 5\", \"46\u00f77=6, 4\", \"63\u00f74=15, 3\", \"96\u00f76=16, 0\", \"19\u00f73=6, 1\"],\n  [\"92\u00f74=23, 0\", \"81\u00f79=9, 0\", \"90\u00f75=18, 0\", \"91\u00f73=30, 1\", \"11\u00f72=5, 1\"],\n  [\"14\u00f73=4, 2\", \"59\u00f76=9, 5\", \"38\u00f77=5, 3\", \"56\u00f78=7, 0\", \"93\u00f75=18, 3\"],\n  [\"88\u00f73=29, 1\", \"57\u00f72=28, 1\", \"15\u00f79=1, 6\", \"37\u00f73=12, 1\", \"57\u00f74=14, 1\"],\n];\nconst populatedRows = [0, 4, 8, 12, 16];\n\nfor (let i = 0; i < populatedRows.length && i < newValues.length; i++) {\n  const rowIndex = populatedRows[i];\n  if (rowIndex >= table.rowCount) break;\n  const rowValues = newValues[i];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(rowIndex, c);\n    const p0 = cell.body.paragraphs.getFirst();\n    p0.insertText(rowValues[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division-answer table cells\n# (in document order) to the new values from the target revision.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date paragraph (first paragraph in the body) ---\n$dateRange = $d.Paragraphs.First.Range\nif ($dateRange.Text -like \"2025-12-18 Thursday*\") {\n    $dateRange.Text = \"2025-12-19 Friday\"\n}\n\n# --- 2. Update the table cells, in row/column (document) order ---\n$tbl = $d.Tables.Item(1)\n\n# New values, listed in the same left-to-right, top-to-bottom order\n# as they appear in the document (only every 4th row is populated).\n$newValues = @(\n    @(\"88\u00f78=11, 0\", \"76\u00f79=8, 4\", \"73\u00f77=10, 3\", \"50\u00f76=8, 2\", \"67\u00f77=9, 4\"),\n    @(\"96\u00f77=13, 5\", \"46\u00f77=6, 4\", \"63\u00f74=15, 3\", \"96\u00f76=16, 0\", \"19\u00f73=6, 1\"),\n    @(\"92\u00f74=23, 0\", \"81\u00f79=9, 0\", \"90\u00f75=18, 0\", \"91\u00f73=30, 1\", \"11\u00f72=5, 1\"),\n    @(\"14\u00f73=4, 2\", \"59\u00f76=9, 5\", \"38\u00f77=5, 3\", \"56\u00f78=7, 0\", \"93\u00f75=18, 3\"),\n    @(\"88\u00f73=29, 1\", \"57\u00f72=28, 1\", \"15\u00f79=1, 6\", \"37\u00f73=12, 1\", \"57\u00f74=14, 1\")\n)\n\n$populatedRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $populatedRows.Count; $i++) {\n    $rowIndex = $populatedRows[$i]\n    $rowValues = $newValues[$i]\n    for ($col = 1; $col -le $rowValues.Count; $col++) {\n        $cell = $tbl.Cell($rowIndex, $col)\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
